# GILDNounTrade.xlsx — record a new trade (row 4) and let the
# "Principle" column (C) re-fit its best-fit width for the longer value.
# traded. Fixed 20 minute trade problem

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New trade row -----------------------------------------------------
# A4: trade close date/time (same date-time display style as A3/G3)
$ws.Range("A4").Value = 42641.540659722225
$ws.Range("A4").NumberFormat = "m/d/yy h:mm"

# B4: Profitable
$ws.Range("B4").Value = $true

# C4: Principle / D4: Start Principle
$ws.Range("C4").Value = 10047.56
$ws.Range("D4").Value = 10023.5

# E4: BuyPrice / F4: SellPrice
$ws.Range("E4").Value = 79.319999999999993
$ws.Range("F4").Value = 78.94

# G4: IsShortSell (same date-time number format slot reused by the sheet)
$ws.Range("G4").Value = $true
$ws.Range("G4").NumberFormat = "m/d/yy h:mm"

# H4: Price Change %
$ws.Range("H4").Value = -0.48

# I4: Strong trade
$ws.Range("I4").Value = $false

# --- Column C best-fit width bump (8.85546875 -> 9) --------------------
# The new "10047.56" value is one character wider than the previous
# longest entry, so the best-fit column width grows by a notch.
$ws.Columns.Item(3).ColumnWidth = 8.14
